# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 5685790.5
$ws.Range("I62").Value = 8336053.5
$ws.Range("J62").Value = 6654.7144
$ws.Range("K62").Value = 8336053.5
$ws.Range("L62").Value = 6654.7144
$ws.Range("M62").Value = -8335429.5
$ws.Range("N62").Value = -7902.7144
# Row 65
$ws.Range("H65").Value = 5685790.5
$ws.Range("I65").Value = 8336053.5
$ws.Range("J65").Value = 6654.7144
$ws.Range("K65").Value = 41680267.5
$ws.Range("L65").Value = 33273.572
$ws.Range("M65").Value = -41677147.5
$ws.Range("N65").Value = -39513.572
# Row 135
$ws.Range("H135").Value = 715550.1
$ws.Range("I135").Value = 870413.4
$ws.Range("J135").Value = 3179
$ws.Range("K135").Value = 7833720.600000001
$ws.Range("L135").Value = 28611
$ws.Range("M135").Value = -7831185.600000001
$ws.Range("N135").Value = -33681
# Row 137
$ws.Range("H137").Value = 420945.4
$ws.Range("I137").Value = 266728.4
$ws.Range("J137").Value = 1006970
$ws.Range("K137").Value = 800185.2000000001
$ws.Range("L137").Value = 3020910
$ws.Range("M137").Value = -797635.2000000001
$ws.Range("N137").Value = -3026010
# Row 138
$ws.Range("H138").Value = 4624.096
$ws.Range("J138").Value = 5949.9424
$ws.Range("L138").Value = 17849.8272
$ws.Range("N138").Value = -28129.8272
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = $null
$ws.Range("N140").Value = 0

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2730.5244
$ws.Range("I32").Value = 1834.0127
$ws.Range("J32").Value = 26338.666
$ws.Range("K32").Value = 1834.0127
$ws.Range("L32").Value = 26338.666
$ws.Range("M32").Value = -1547.0127
$ws.Range("N32").Value = -26912.666
# Row 46
$ws.Range("H46").Value = 6977.4287
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 7807
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 7807
$ws.Range("M46").Value = -1681
$ws.Range("N46").Value = -8445
# Row 61
$ws.Range("H61").Value = 1966.0834
$ws.Range("I61").Value = 1719.7805
$ws.Range("J61").Value = 3408.7144
$ws.Range("K61").Value = 1719.7805
$ws.Range("L61").Value = 3408.7144
$ws.Range("M61").Value = -1507.7805
$ws.Range("N61").Value = -3832.7144
# Row 136
$ws.Range("H136").Value = 1966.0834
$ws.Range("I136").Value = 1719.7805
$ws.Range("J136").Value = 3408.7144
$ws.Range("K136").Value = 5159.3415
$ws.Range("L136").Value = 10226.1432
$ws.Range("M136").Value = -2609.3415
$ws.Range("N136").Value = -15326.1432

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1999.875
$ws.Range("I105").Value = 1125
$ws.Range("J105").Value = 2874.75
$ws.Range("K105").Value = 1125
$ws.Range("L105").Value = 2874.75
$ws.Range("M105").Value = 622
$ws.Range("N105").Value = -6368.75
# Row 134
$ws.Range("H134").Value = 48906.39
$ws.Range("I134").Value = 4991.9473
$ws.Range("J134").Value = 257500
$ws.Range("K134").Value = 14975.8419
$ws.Range("L134").Value = 772500
$ws.Range("M134").Value = -12440.8419
$ws.Range("N134").Value = -777570

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 53144.75
$ws.Range("I31").Value = 1555.0625
$ws.Range("J31").Value = 259503.5
$ws.Range("K31").Value = 1555.0625
$ws.Range("L31").Value = 259503.5
$ws.Range("M31").Value = -1260.0625
$ws.Range("N31").Value = -260093.5
# Row 34
$ws.Range("H34").Value = 53144.75
$ws.Range("I34").Value = 1555.0625
$ws.Range("J34").Value = 259503.5
$ws.Range("K34").Value = 1555.0625
$ws.Range("L34").Value = 259503.5
$ws.Range("M34").Value = -1353.0625
$ws.Range("N34").Value = -259907.5
# Row 35
$ws.Range("H35").Value = 4467.727
$ws.Range("I35").Value = 1717.8572
$ws.Range("K35").Value = 1717.8572
$ws.Range("M35").Value = -1423.8572
# Row 58
$ws.Range("H58").Value = 135798.2
$ws.Range("I58").Value = 171298.6
$ws.Range("J58").Value = 4890.5
$ws.Range("K58").Value = 171298.6
$ws.Range("L58").Value = 4890.5
$ws.Range("M58").Value = -171095.6
$ws.Range("N58").Value = -5296.5
# Row 132
$ws.Range("H132").Value = 2530.8772
$ws.Range("I132").Value = 2221.275
$ws.Range("K132").Value = 6663.825000000001
$ws.Range("M132").Value = -4133.825000000001
# Row 134
$ws.Range("H134").Value = 366028.8
$ws.Range("I134").Value = 210753.83
$ws.Range("J134").Value = 1430771.4
$ws.Range("K134").Value = 632261.49
$ws.Range("L134").Value = 4292314.199999999
$ws.Range("M134").Value = -629726.49
$ws.Range("N134").Value = -4297384.199999999
# Row 136
$ws.Range("H136").Value = 135798.2
$ws.Range("I136").Value = 171298.6
$ws.Range("J136").Value = 4890.5
$ws.Range("K136").Value = 513895.8
$ws.Range("L136").Value = 14671.5
$ws.Range("M136").Value = -511345.8
$ws.Range("N136").Value = -19771.5

$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 10035001
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 10035001
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = $null
$ws.Range("M9").Value = 30105003
$ws.Range("N9").Value = -30105451
# Row 22
$ws.Range("H22").Value = 1861.625
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 2070.4285
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 6211.2855
$ws.Range("M22").Value = -1031
$ws.Range("N22").Value = -6549.2855
# Row 27
$ws.Range("H27").Value = 1861.625
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 2070.4285
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 6211.2855
$ws.Range("M27").Value = -1098
$ws.Range("N27").Value = -6415.2855
# Row 99
$ws.Range("H99").Value = 3245
$ws.Range("I99").Value = 1830
$ws.Range("J99").Value = 4660
$ws.Range("K99").Value = 5490
$ws.Range("L99").Value = 13980
$ws.Range("M99").Value = -3244
$ws.Range("N99").Value = -18472

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 231272.55
$ws.Range("I132").Value = 229433.05
$ws.Range("J132").Value = 251507
$ws.Range("K132").Value = 688299.1499999999
$ws.Range("L132").Value = 754521
$ws.Range("M132").Value = -685769.1499999999
$ws.Range("N132").Value = -759581

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 804.5
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 835.4
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 835.4
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -1425.4
# Row 27
$ws.Range("H27").Value = 804.5
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 835.4
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 835.4
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -1049.4
# Row 136
$ws.Range("H136").Value = 309852.34
$ws.Range("I136").Value = 346582.78
$ws.Range("J136").Value = 5514.5713
$ws.Range("K136").Value = 1039748.34
$ws.Range("L136").Value = 16543.7139
$ws.Range("M136").Value = -1037198.34
$ws.Range("N136").Value = -21643.7139

$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = $null
$ws.Range("N109").Value = 0
# Row 122
$ws.Range("H122").Value = 19610024
$ws.Range("I122").Value = 27028540
$ws.Range("J122").Value = 3944.4285
$ws.Range("K122").Value = 81085620
$ws.Range("L122").Value = 11833.2855
$ws.Range("M122").Value = -81083170
$ws.Range("N122").Value = -16733.2855
